$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing rows (columns C, D, E, F) ---

# Row 437
$ws.Range("C437").Value = 106
$ws.Range("D437").Value = 87.95999999999999
$ws.Range("E437").Value = 220.6
$ws.Range("F437").Value = 183.06

# Row 440
$ws.Range("C440").Value = 100.19
$ws.Range("D440").Value = 84.56999999999999
$ws.Range("E440").Value = 226.53
$ws.Range("F440").Value = 191.21

# Row 443
$ws.Range("C443").Value = 102.6
$ws.Range("D443").Value = 85.52
$ws.Range("E443").Value = 234.19
$ws.Range("F443").Value = 195.19

# Row 444
$ws.Range("C444").Value = 100.91
$ws.Range("D444").Value = 83.54000000000001
$ws.Range("E444").Value = 237.75
$ws.Range("F444").Value = 196.82

# Row 445
$ws.Range("D445").Value = 81.78
$ws.Range("E445").Value = 243.08
$ws.Range("F445").Value = 200.66

# Row 446
$ws.Range("F446").Value = 204.34

# Row 449
$ws.Range("C449").Value = 96.84999999999999
$ws.Range("D449").Value = 81.62
$ws.Range("E449").Value = 250.88
$ws.Range("F449").Value = 211.44

# Row 450
$ws.Range("F450").Value = 216.96

# Row 451
$ws.Range("C451").Value = 101.96
$ws.Range("D451").Value = 86.47
$ws.Range("E451").Value = 258.35
$ws.Range("F451").Value = 219.1

# Row 452
$ws.Range("C452").Value = 104.26
$ws.Range("D452").Value = 88.83
$ws.Range("E452").Value = 257.83
$ws.Range("F452").Value = 219.67

# --- Complete row 453 (previously only had A453, B453) ---
$ws.Range("C453").Value = 108.06
$ws.Range("D453").Value = 92.42
$ws.Range("E453").Value = 258.18
$ws.Range("F453").Value = 220.81

# --- New row 454 with a new month ---
$ws.Range("A454").Value = "01-09-2021"
$ws.Range("B454").Value = 783.63
